# "Updated for 4 th apr"
# Duplicate the "3-Apr" report sheet to create a new "4-Apr" sheet (same
# layout/data as 3-Apr had), append it as the last tab, and restore the
# view/selection state on both sheets to match what Excel saved:
#   - "3-Apr"  (no longer the active tab) -> selection moves to G20
#   - "4-Apr"  (new active/selected tab)  -> selection at F14

$wb = $excel.ActiveWorkbook

$sourceSheet = $wb.Worksheets.Item("3-Apr")

# Copy the sheet to right after itself -> lands as the new last tab.
$sourceSheet.Copy($null, $sourceSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "4-Apr"

# "3-Apr" keeps its data untouched, only its saved selection changes.
$sourceSheet.Activate()
$sourceSheet.Range("G20").Select() | Out-Null

# The new "4-Apr" sheet becomes the active/selected tab, cursor on F14.
$newSheet.Activate()
$newSheet.Range("F14").Select() | Out-Null
